# Update column F (dSF) values on the active sheet to reflect the
# re-pulled data / mean calculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 2
    4  = -4
    5  = -5
    7  = 3
    8  = -2
    10 = 3
    11 = -6
    12 = 1
    14 = -2
    15 = 4
    16 = 2
    17 = 3
    21 = 1
    22 = -1
    23 = 6
    24 = -2
    25 = 3
    26 = 1
    27 = 3
    28 = -3
    29 = 2
    30 = 4
    31 = 5
    32 = -5
    33 = 3
    34 = 3
    35 = 1
    36 = 2
    37 = 2
    38 = 5
    39 = 3
    40 = -1
    41 = -1
    42 = -5
    43 = 1
    45 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
